$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.855461368907875
$ws.Range("C2").Value = 0.05571496487019689

$ws.Range("B3").Value = 5.373821717978156
$ws.Range("C3").Value = 0.0509700646315515

$ws.Range("B4").Value = 5.39224651444028
$ws.Range("C4").Value = 0.04786369706241319

$ws.Range("B5").Value = 5.242322013566242
$ws.Range("C5").Value = 0.04518225044339162

$ws.Range("B6").Value = 5.250397634109723
$ws.Range("C6").Value = 0.04377345463725519

$ws.Range("B7").Value = 5.209812860701463
$ws.Range("C7").Value = 0.04182734253918256

$ws.Range("B8").Value = 4.823498661210245
$ws.Range("C8").Value = 0.04071928958785433

$ws.Range("B9").Value = 5.06411657742001
$ws.Range("C9").Value = 0.03963256045586165

$ws.Range("B10").Value = 4.852897753594312
$ws.Range("C10").Value = 0.03897010148013704
